$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preventing Excel from
# auto-coercing numeric-looking strings (e.g. "1.00") into numbers,
# and without leaving a lasting number-format style on the cell.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '69.562.56'
Set-TextValue $ws.Range('E2') '  -1.00%  '
Set-TextValue $ws.Range('D3') '3.507.29'
Set-TextValue $ws.Range('E3') '  -2.59%  '
Set-TextValue $ws.Range('E4') '  -0.07%  '
Set-TextValue $ws.Range('D5') '585.12'
Set-TextValue $ws.Range('E5') '  +0.41%  '
Set-TextValue $ws.Range('D6') '184.88'
Set-TextValue $ws.Range('E6') '  -2.62%  '
Set-TextValue $ws.Range('D7') '3.493.20'
Set-TextValue $ws.Range('E7') '  -2.84%  '
Set-TextValue $ws.Range('E8') '  -3.22%  '
Set-TextValue $ws.Range('D10') '0.201'
Set-TextValue $ws.Range('E10') '  +10.54%  '
Set-TextValue $ws.Range('D11') '0.646'
Set-TextValue $ws.Range('E11') '  -2.96%  '
Set-TextValue $ws.Range('D12') '53.83'
Set-TextValue $ws.Range('E12') '  -4.02%  '
Set-TextValue $ws.Range('D13') '0.0000307'
Set-TextValue $ws.Range('E13') '  -1.36%  '
Set-TextValue $ws.Range('D14') '9.43'
Set-TextValue $ws.Range('E14') '  -3.03%  '
Set-TextValue $ws.Range('D15') '4.054.40'
Set-TextValue $ws.Range('E15') '  -3.10%  '
Set-TextValue $ws.Range('D16') '19.29'
Set-TextValue $ws.Range('E16') '  -3.69%  '
Set-TextValue $ws.Range('D17') '69.469.99'
Set-TextValue $ws.Range('E17') '  -1.13%  '
Set-TextValue $ws.Range('D18') '3.493.99'
Set-TextValue $ws.Range('E18') '  -3.05%  '
Set-TextValue $ws.Range('D19') '12.32'
Set-TextValue $ws.Range('E19') '  -3.10%  '
Set-TextValue $ws.Range('E20') '  -1.34%  '
Set-TextValue $ws.Range('D21') '543.65'
Set-TextValue $ws.Range('E21') '  +11.27%  '
Set-TextValue $ws.Range('E22') '  -3.83%  '
Set-TextValue $ws.Range('D23') '18.47'
Set-TextValue $ws.Range('E23') '  -7.78%  '
Set-TextValue $ws.Range('D24') '4.53'
Set-TextValue $ws.Range('E24') '  +2.99%  '
Set-TextValue $ws.Range('D25') '4.87'
Set-TextValue $ws.Range('E25') '  -1.99%  '
Set-TextValue $ws.Range('D26') '94.97'
Set-TextValue $ws.Range('E26') '  -2.35%  '
Set-TextValue $ws.Range('D27') '2.98'
Set-TextValue $ws.Range('E27') '  -0.32%  '
Set-TextValue $ws.Range('D28') '11.06'
Set-TextValue $ws.Range('E28') '  -0.62%  '
Set-TextValue $ws.Range('D29') '9.11'
Set-TextValue $ws.Range('E29') '  -4.42%  '
Set-TextValue $ws.Range('D30') '31.99'
Set-TextValue $ws.Range('E30') '  -1.35%  '
Set-TextValue $ws.Range('D31') '7.30'
Set-TextValue $ws.Range('E31') '  -4.16%  '
Set-TextValue $ws.Range('D32') '12.60'
Set-TextValue $ws.Range('E32') '  +2.33%  '
Set-TextValue $ws.Range('D33') '64.00'
Set-TextValue $ws.Range('E33') '  -3.67%  '
Set-TextValue $ws.Range('E34') '  -4.86%  '
Set-TextValue $ws.Range('D35') '532.05'
Set-TextValue $ws.Range('E35') '  -7.89%  '
Set-TextValue $ws.Range('D36') '0.407'
Set-TextValue $ws.Range('E36') '  +1.97%  '
Set-TextValue $ws.Range('B37') 'InjectiveProtocol'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D37') '38.32'
Set-TextValue $ws.Range('E37') '  -1.91%  '
Set-TextValue $ws.Range('E38') '  -0.13%  '
Set-TextValue $ws.Range('B39') 'Fetch.AI'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D39') '3.05'
Set-TextValue $ws.Range('E39') '  +4.87%  '
Set-TextValue $ws.Range('D40') '0.0₃0767'
Set-TextValue $ws.Range('E40') '  -6.25%  '
Set-TextValue $ws.Range('E41') '  -3.56%  '
Set-TextValue $ws.Range('B42') 'Kaspa'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D42') '0.134'
Set-TextValue $ws.Range('E42') '  -1.85%  '
Set-TextValue $ws.Range('B43') 'Stacks'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D43') '3.39'
Set-TextValue $ws.Range('E43') '  -2.14%  '
Set-TextValue $ws.Range('B44') 'Maker'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D44') '3.339.52'
Set-TextValue $ws.Range('E44') '  +3.56%  '
Set-TextValue $ws.Range('D45') '3.54'
Set-TextValue $ws.Range('E45') '  +5.92%  '
Set-TextValue $ws.Range('D46') '2.98'
Set-TextValue $ws.Range('E46') '  -3.26%  '
Set-TextValue $ws.Range('D47') '0.0441'
Set-TextValue $ws.Range('E47') '  -1.54%  '
Set-TextValue $ws.Range('D48') '9.01'
Set-TextValue $ws.Range('E48') '  -6.20%  '
Set-TextValue $ws.Range('E49') '  -3.28%  '
Set-TextValue $ws.Range('D50') '1.00'
Set-TextValue $ws.Range('E50') '  +0.11%  '
Set-TextValue $ws.Range('D51') '136.51'
Set-TextValue $ws.Range('E51') '  +0.15%  '
